$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06670480677062823
$ws.Range("H2").Value = 3.016451874312913
$ws.Range("I2").Value = -336.8347316467106
$ws.Range("G3").Value = 0.07330992072963301
$ws.Range("H3").Value = 7.170300594181287
$ws.Range("G4").Value = -0.03294016938801222
$ws.Range("H4").Value = 26.62521564245726
$ws.Range("G5").Value = -0.03129202290480846
$ws.Range("H5").Value = -15.58881843488774
$ws.Range("G6").Value = -0.09396381958128665
$ws.Range("H6").Value = 11.38932692166874
$ws.Range("G7").Value = -0.08265759590487787
$ws.Range("H7").Value = 9.53645685219526
$ws.Range("G8").Value = -0.3476890560383758
$ws.Range("H8").Value = 5.237186200789197
$ws.Range("G9").Value = -0.3769940312045026
$ws.Range("H9").Value = 3.357345533586839
$ws.Range("G10").Value = 0.02876582427721338
$ws.Range("H10").Value = 42.4685613899169
$ws.Range("G11").Value = 0.0447836922367937
$ws.Range("H11").Value = 97.33300803036923
$ws.Range("G12").Value = 0.2204289536178456
$ws.Range("H12").Value = -0.5946389650944507
$ws.Range("G13").Value = 0.2309679568309875
$ws.Range("H13").Value = 2.558162283729016
$ws.Range("G14").Value = -0.06376907605958552
$ws.Range("H14").Value = -51.45213276741751
$ws.Range("G15").Value = -0.03845660957441992
$ws.Range("H15").Value = 19.37191687799868
$ws.Range("G16").Value = 0.2157936062243941
$ws.Range("H16").Value = 1.514868497466223
$ws.Range("G17").Value = 0.2143220118156841
$ws.Range("H17").Value = -2.822947199916408
$ws.Range("G18").Value = 0.07517709543238489
$ws.Range("H18").Value = 2.953776683646371
$ws.Range("G19").Value = 0.0764080057458804
$ws.Range("H19").Value = 1.42025641196491
$ws.Range("G20").Value = -0.0811046754097196
$ws.Range("H20").Value = -8.1655614022296
$ws.Range("G21").Value = -0.08626032590905523
$ws.Range("H21").Value = 0.3563701529764416
$ws.Range("G22").Value = 0.07397295211361647
$ws.Range("H22").Value = 0.6408749974718129
$ws.Range("G23").Value = 0.06739025330740724
$ws.Range("H23").Value = -1.377206083962102
$ws.Range("G24").Value = 0.06403731134768341
$ws.Range("H24").Value = -3.865391778840805
$ws.Range("G25").Value = 0.06848276779957174
$ws.Range("H25").Value = 25.01405688049451
$ws.Range("G26").Value = 0.1158078167661229
$ws.Range("H26").Value = -2.967110272928283
$ws.Range("G27").Value = 0.1210678613475955
$ws.Range("H27").Value = 6.330288468543478
$ws.Range("G28").Value = 0.1389334232484681
$ws.Range("H28").Value = 7.489200226186403
$ws.Range("G29").Value = 0.1443895941216611
$ws.Range("H29").Value = -4.276585341513358
$ws.Range("G30").Value = 0.08428650176641249
$ws.Range("H30").Value = -0.02467776273938008
$ws.Range("G31").Value = 0.08428650176641249
$ws.Range("H31").Value = 3.180951443921745
$ws.Range("G32").Value = 0.05377260711433124
$ws.Range("H32").Value = 0.7734460567699114
$ws.Range("G33").Value = 0.05644676705635598
$ws.Range("H33").Value = 2.179204903091989
$ws.Range("G34").Value = 0.02410073199091677
$ws.Range("H34").Value = 38.8511438349663
$ws.Range("G35").Value = 0.02410073199091677
$ws.Range("H35").Value = 42.60207511346352
$ws.Range("G36").Value = -0.02360738497116695
$ws.Range("H36").Value = 18.72330649821934
$ws.Range("G37").Value = -0.02711740583071
$ws.Range("H37").Value = 18.4779198589934
$ws.Range("G38").Value = 0.08037878804760902
$ws.Range("H38").Value = 2.67755767027846
$ws.Range("G39").Value = 0.08417697981510645
$ws.Range("H39").Value = 8.272992547135868
$ws.Range("G40").Value = 0.06190496785613346
$ws.Range("H40").Value = -6.500720016109855
$ws.Range("G41").Value = 0.07339705216844888
$ws.Range("H41").Value = 12.87887884232451
$ws.Range("G42").Value = 0.08300680658778022
$ws.Range("H42").Value = 6.704495130924405
$ws.Range("G43").Value = 0.07522735866935555
$ws.Range("H43").Value = -6.158753396946476
$ws.Range("G44").Value = 0.08417006730378612
$ws.Range("H44").Value = -4.620237324517281
$ws.Range("G45").Value = 0.09404507662440652
$ws.Range("H45").Value = 4.051149580269861
$ws.Range("G46").Value = -0.002060324764616932
$ws.Range("H46").Value = 24.70182635872525
$ws.Range("G47").Value = -0.009230174638582202
$ws.Range("H47").Value = -19220.74199795027
$ws.Range("G48").Value = -0.1070881896171063
$ws.Range("H48").Value = -11.4281555281762
$ws.Range("G49").Value = -0.09922648753949648
$ws.Range("H49").Value = 9.445221926896503
$ws.Range("G50").Value = 0.1663052402578979
$ws.Range("H50").Value = -2.459864005318104
$ws.Range("G51").Value = 0.1749837019877335
$ws.Range("H51").Value = 3.031670354912904
$ws.Range("G52").Value = 0.0714218946801269
$ws.Range("H52").Value = 0.648249537585643
$ws.Range("G53").Value = 0.06649836295586489
$ws.Range("H53").Value = 3.393014480306258
$ws.Range("G54").Value = -0.1208743071529234
$ws.Range("H54").Value = 5.421278313618873
$ws.Range("G55").Value = -0.1208743071529234
$ws.Range("H55").Value = -3.780454558428616
$ws.Range("G56").Value = 0.1928873954466858
$ws.Range("H56").Value = 1.508550533292732
$ws.Range("G57").Value = 0.2044254390504687
$ws.Range("H57").Value = 2.773723067853299

Write-Output "Updated G/H/I columns for rows 2-57"
